$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Insert a new column before column B for the "author" field,
# shifting date/yoast_metadesc/excerpt/category/tags/url_path right by one.
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ClearFormats()

$ws.Range("B1").Value = "author"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108

$authors = @{
    2 = "Exadel Digital Transformation Team"
    3 = "Sergey Derugo"
    4 = "Maryna Shantur"
    5 = "Nikita Basalaev"
    6 = "Nikita Basalaev"
    7 = "Exadel Digital Transformation Team"
    8 = "Alexey Chumakov"
    9 = "Iryna Ason"
    10 = "Exadel Digital Transformation Team"
    11 = "Exadel Digital Transformation Team"
    12 = "Exadel Digital Marketing Technology Team"
    13 = "Liubou Masiuk"
    14 = "Exadel Digital Marketing Technology Team"
    15 = "Jonathan Fries"
    16 = "Exadel Digital Transformation Team"
    17 = "Iryna Ason"
    18 = "Jonathan Fries"
    19 = "Jonathan Fries"
    20 = "Exadel Digital Transformation Team"
    21 = "Lev Shur"
    22 = "Lev Shur"
    23 = "Nikita Basalaev"
    24 = "Alexander Bestsenny"
    25 = "Nikita Basalaev"
    26 = "Lev Shur"
    27 = "Nikita Basalaev"
    28 = "Nikita Basalaev"
    29 = "Nikita Basalaev"
    30 = "Nikita Basalaev"
    31 = "Nikita Basalaev"
    32 = "Nikita Basalaev"
    33 = "Olga Zhuravleva"
    34 = "Exadel Innovations Team"
    35 = "Exadel Innovations Team"
    36 = "Denis Glushkov"
    37 = "Exadel Innovations Team"
    38 = "Exadel Innovations Team"
    39 = "Exadel Innovations Team"
    40 = "Dmitry Binunsky"
}

foreach ($row in $authors.Keys) {
    $ws.Cells.Item($row, 2).Value = $authors[$row]
}

